$wb = $excel.ActiveWorkbook

# Update the "Status" text for both rows from "Ready for handoff" to "In Translation"
# across the Overview sheet (per-language status columns) and the two
# per-language detail sheets (Status column).

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# The shortened status text narrows the autosized "Status" columns.
$overview.Range("E1").EntireColumn.ColumnWidth = 12.5
$overview.Range("F1").EntireColumn.ColumnWidth = 12.5
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
